# Update Thailand Premier League odds bases (15-06-2024 21:10).
# Rows 234 and 236, and rows 235 and 237, swap their full set of match
# data (everything except the running index in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

Swap-RowData $ws 234 236
Swap-RowData $ws 235 237
